$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 80, shifting existing rows 80:155 down to 81:156
$ws.Rows("80:80").Insert()

# Populate the newly inserted row 80 with the new data point
$ws.Range("A80").Value = 5
$ws.Range("B80").Value = "Macroferia Regional de Talca"
$ws.Range("C80").Value = "Maule"
$ws.Range("D80").Value = 45271
$ws.Range("E80").Value = 7
$ws.Range("F80").Value = "Fruta"
$ws.Range("G80").Value = 100101
$ws.Range("H80").Value = "Berries"
$ws.Range("I80").Value = 100101001
$ws.Range("J80").Value = "Arándano (blue)"
$ws.Range("K80").Value = "Sin especificar"
$ws.Range("L80").Value = "Primera"
$ws.Range("M80").Value = 100
$ws.Range("N80").Value = 4000
$ws.Range("O80").Value = 4000
$ws.Range("P80").Value = 4000
$ws.Range("Q80").Value = "$/bandeja 2 kilos"
$ws.Range("R80").Value = "Provincia de Curicó"
$ws.Range("S80").Value = 2000
$ws.Range("T80").Value = 2
